$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells requiring a style/type change: copy formatting from a stable same-style cell, then set value ---
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("I22").Copy($ws.Range("D22"))
$ws.Range("L22").Copy($ws.Range("E22"))
$ws.Range("I22").Copy($ws.Range("G22"))
$ws.Range("L22").Copy($ws.Range("H22"))
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# --- Set final values for all changed cells ---
$ws.Range("N14").Value = -50
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -36.666666666666
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("I16").Value = 188
$ws.Range("J16").Value = 108
$ws.Range("K16").Value = 74.074074074074
$ws.Range("L16").Value = 54.098360655737
$ws.Range("M16").Value = -25.396825396825
$ws.Range("N16").Value = -82.462686567164
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 80
$ws.Range("I17").Value = 344
$ws.Range("J17").Value = 244
$ws.Range("K17").Value = 40.983606557377
$ws.Range("L17").Value = 50.218340611353
$ws.Range("M17").Value = 62.264150943396
$ws.Range("N17").Value = 2.686567164179
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 141
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 42.424242424242
$ws.Range("L18").Value = 0.714285714285
$ws.Range("M18").Value = -56.748466257668
$ws.Range("N18").Value = -89.350453172205
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -15.789473684210
$ws.Range("I19").Value = 344
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 37.6
$ws.Range("L19").Value = 34.901960784313
$ws.Range("M19").Value = 5.846153846153
$ws.Range("N19").Value = -34.476190476190
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 240
$ws.Range("J20").Value = 219
$ws.Range("K20").Value = 9.589041095890
$ws.Range("L20").Value = 71.428571428571
$ws.Range("M20").Value = -10.112359550561
$ws.Range("N20").Value = -92.412266835283
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = 5.882352941176
$ws.Range("I21").Value = 1281
$ws.Range("J21").Value = 945
$ws.Range("K21").Value = 35.555555555555
$ws.Range("L21").Value = 40.769230769230
$ws.Range("M21").Value = -8.434596140100
$ws.Range("N21").Value = -80.167208546214
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = 62.5
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 91.666666666666
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 22.222222222222
$ws.Range("I24").Value = 1504
$ws.Range("J24").Value = 931
$ws.Range("K24").Value = 61.546723952739
$ws.Range("L24").Value = 72.279495990836
$ws.Range("M24").Value = 107.734806629834
$ws.Range("C25").Value = 19
$ws.Range("E25").Value = 72.727272727272
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = 6.122448979591
$ws.Range("I25").Value = 507
$ws.Range("J25").Value = 468
$ws.Range("K25").Value = 8.333333333333
$ws.Range("L25").Value = 26.433915211970
$ws.Range("M25").Value = -11.363636363636
$ws.Range("C26").Value = "0"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("L26").Value = 6.451612903225
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 49
$ws.Range("J27").Value = 66
$ws.Range("K27").Value = -25.757575757575
$ws.Range("L27").Value = 4.255319148936
$ws.Range("N28").Value = -35.483870967741
$ws.Range("N29").Value = -61.538461538461
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
